$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Version and Date values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.2.0-ballot"
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# --- Include ValueSet #0: append version to the ValueSet URL ---
$inc0 = $wb.Worksheets.Item("Include ValueSet #0")
$inc0.Range("A2").Value = "https://hl7.fr/ig/fhir/core/ValueSet/fr-core-vs-encounter-type|2.1.0"

# --- Include ValueSet #2: append version to the ValueSet URL ---
$inc2 = $wb.Worksheets.Item("Include ValueSet #2")
$inc2.Range("A2").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-type-evenement-ssiad-cisis|20250624152100"
